$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Load Step Number row (row 2): B2 1 -> 2, C2 2 -> 3
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = 3

# Update Bending [in-lbf] formula (row 4): =975*12 -> =976*12
$ws.Range("B4").Formula = "=976*12"

# Update the active selection to B6 (matches saved view state in the file)
$ws.Range("B6").Select()
